# Updates cryptos list: refresh Price (column D) and Volume(1h) (column E)
# values for each coin row, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are stored as plain text (e.g. "0.590", "67.442.35")
# rather than numbers, so force Text format before writing to avoid Excel
# auto-converting them to numeric values and stripping significant trailing
# zeros / thousands-dot formatting.
$priceCells = @('D2', 'D3', 'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D28', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50', 'D51')
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.872.96'
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('D3').Value = '3.325.48'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('D4').Value = '0.996'
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = '587.21'
$ws.Range('E5').Value = '  +5.80%  '
$ws.Range('D6').Value = '182.91'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('D8').Value = '0.590'
$ws.Range('E8').Value = '  +3.27%  '
$ws.Range('D9').Value = '3.309.69'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').Value = '0.580'
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '46.54'
$ws.Range('E12').Value = '  +2.46%  '
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  +7.04%  '
$ws.Range('D14').Value = '638.91'
$ws.Range('E14').Value = '  +11.77%  '
$ws.Range('D15').Value = '3.821.21'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').Value = '8.46'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '67.803.25'
$ws.Range('E17').Value = '  +3.25%  '
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').Value = '3.298.74'
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = '17.72'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').Value = '10.97'
$ws.Range('E21').Value = '  +2.43%  '
$ws.Range('D22').Value = '0.901'
$ws.Range('E22').Value = '  +1.95%  '
$ws.Range('D23').Value = '17.75'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '5.04'
$ws.Range('E24').Value = '  +1.54%  '
$ws.Range('D25').Value = '98.11'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').Value = '4.00'
$ws.Range('E26').Value = '  +2.07%  '
$ws.Range('D27').Value = '2.81'
$ws.Range('E27').Value = '  +5.99%  '
$ws.Range('D28').Value = '9.66'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('D29').Value = '32.95'
$ws.Range('E29').Value = '  +8.80%  '
$ws.Range('D30').Value = '8.59'
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('D31').Value = '6.70'
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').Value = '600.61'
$ws.Range('E32').Value = '  +8.19%  '
$ws.Range('D33').Value = '3.940.57'
$ws.Range('E33').Value = '  +5.00%  '
$ws.Range('D34').Value = '3.63'
$ws.Range('E34').Value = '  -0.68%  '
$ws.Range('D35').Value = '10.97'
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('D37').Value = '0.997'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').Value = '55.77'
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').Value = '3.30'
$ws.Range('E39').Value = '  +6.57%  '
$ws.Range('D40').Value = '2.72'
$ws.Range('E40').Value = '  +7.06%  '
$ws.Range('D41').Value = '0.129'
$ws.Range('E41').Value = '  +2.73%  '
$ws.Range('D42').Value = '33.07'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('D43').Value = '0.0₃0691'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('D44').Value = '3.35'
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').Value = '0.338'
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('D46').Value = '0.0414'
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('D47').Value = '0.129'
$ws.Range('E47').Value = '  +2.37%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = '2.55'
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('D50').Value = '1.34'
$ws.Range('E50').Value = '  +9.16%  '
$ws.Range('D51').Value = '131.07'
$ws.Range('E51').Value = '  +5.07%  '
